# Applies the EffectData.xlsx edit:
#  - Zero out columns AB:AF for rows 11..70 (was 50, now 0)
#  - Split column AG (33) off from AF's width band and widen it
#  - Update the frozen-pane scroll position / active selection on sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# --- Data edit: AB11:AF70 -> 0 (was 50) -------------------------------
$ws.Range("AB11:AF70").Value = 0

# --- Column width: give AG (col 33) its own (wider) width -------------
$ws.Columns.Item(33).ColumnWidth = 15

# --- View state: scroll the unfrozen pane to column Y, row 11, then
#     select AB11:AB70 (matches the sheetView/selection captured on save)
$ws.Activate()
$window = $excel.ActiveWindow
$window.FreezePanes = $false
$ws.Range("Y11").Select()
$window.FreezePanes = $true
$ws.Range("AB11:AB70").Select()
